# Apply the "Updated symbol list" refresh: updated prices/volumes/hour stamps,
# and the Coin/Link rows 16-20 shifted down one position with a new TigerCash entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name / Link URL) -- Excel keeps these as text natively ---
$plainUpdates = @(
    @{ Ref = 'B16'; Value = 'TigerCash' },
    @{ Ref = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Ref = 'B17'; Value = 'UpBots' },
    @{ Ref = 'C17'; Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt' },
    @{ Ref = 'B18'; Value = 'LEO' },
    @{ Ref = 'C18'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Ref = 'B19'; Value = 'BTSEToken' },
    @{ Ref = 'C19'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' },
    @{ Ref = 'B20'; Value = 'One' },
    @{ Ref = 'C20'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
)
foreach ($u in $plainUpdates) {
    $ws.Range($u.Ref).Value = $u.Value
}

# --- Numeric-looking text cells (Price / Volume% / Hora) ---
# These are stored as text in the sheet (e.g. "246.63", "0.56%", "18"), so force
# the cell to Text format before writing, then restore the default style so no
# stray formatting is left behind on the cell.
$textUpdates = @(
    @{ Ref = 'D2'; Value = '246.63' },
    @{ Ref = 'E2'; Value = '0.56%' },
    @{ Ref = 'G2'; Value = '18' },
    @{ Ref = 'D3'; Value = '29.84' },
    @{ Ref = 'E3'; Value = '10.74%' },
    @{ Ref = 'G3'; Value = '18' },
    @{ Ref = 'D4'; Value = '5.168' },
    @{ Ref = 'E4'; Value = '1.90%' },
    @{ Ref = 'G4'; Value = '18' },
    @{ Ref = 'D5'; Value = '0.05712' },
    @{ Ref = 'E5'; Value = '0.36%' },
    @{ Ref = 'G5'; Value = '18' },
    @{ Ref = 'D6'; Value = '6.614' },
    @{ Ref = 'E6'; Value = '1.93%' },
    @{ Ref = 'G6'; Value = '18' },
    @{ Ref = 'D7'; Value = '3.076' },
    @{ Ref = 'E7'; Value = '2.36%' },
    @{ Ref = 'G7'; Value = '18' },
    @{ Ref = 'D8'; Value = '0.8562' },
    @{ Ref = 'E8'; Value = '4.47%' },
    @{ Ref = 'G8'; Value = '18' },
    @{ Ref = 'D9'; Value = '0.8693' },
    @{ Ref = 'E9'; Value = '3.07%' },
    @{ Ref = 'G9'; Value = '18' },
    @{ Ref = 'D10'; Value = '0.1366' },
    @{ Ref = 'E10'; Value = '2.86%' },
    @{ Ref = 'G10'; Value = '18' },
    @{ Ref = 'D11'; Value = '0.07074' },
    @{ Ref = 'E11'; Value = '2.91%' },
    @{ Ref = 'G11'; Value = '18' },
    @{ Ref = 'D12'; Value = '0.02924' },
    @{ Ref = 'E12'; Value = '3.24%' },
    @{ Ref = 'G12'; Value = '18' },
    @{ Ref = 'D13'; Value = '0.09391' },
    @{ Ref = 'E13'; Value = '-0.04%' },
    @{ Ref = 'G13'; Value = '18' },
    @{ Ref = 'D14'; Value = '0.001528' },
    @{ Ref = 'E14'; Value = '0.26%' },
    @{ Ref = 'G14'; Value = '18' },
    @{ Ref = 'D15'; Value = '0.04179' },
    @{ Ref = 'E15'; Value = '2.40%' },
    @{ Ref = 'G15'; Value = '18' },
    @{ Ref = 'D16'; Value = '0.006092' },
    @{ Ref = 'E16'; Value = '0.00%' },
    @{ Ref = 'G16'; Value = '18' },
    @{ Ref = 'D17'; Value = '0.007489' },
    @{ Ref = 'E17'; Value = '3,766.06%' },
    @{ Ref = 'G17'; Value = '18' },
    @{ Ref = 'D18'; Value = '3.487' },
    @{ Ref = 'E18'; Value = '-0.67%' },
    @{ Ref = 'G18'; Value = '18' },
    @{ Ref = 'D19'; Value = '2.268' },
    @{ Ref = 'E19'; Value = '1.72%' },
    @{ Ref = 'G19'; Value = '18' },
    @{ Ref = 'D20'; Value = '0.0005990' },
    @{ Ref = 'E20'; Value = '-0.63%' },
    @{ Ref = 'G20'; Value = '18' },
    @{ Ref = 'G21'; Value = '18' },
    @{ Ref = 'G22'; Value = '18' },
    @{ Ref = 'D23'; Value = '0.1300' },
    @{ Ref = 'E23'; Value = '0.22%' },
    @{ Ref = 'G23'; Value = '18' },
    @{ Ref = 'D24'; Value = '3.469' },
    @{ Ref = 'E24'; Value = '-2.41%' },
    @{ Ref = 'G24'; Value = '18' },
    @{ Ref = 'E25'; Value = '0.49%' },
    @{ Ref = 'G25'; Value = '18' },
    @{ Ref = 'D26'; Value = '0.005023' },
    @{ Ref = 'E26'; Value = '26.67%' },
    @{ Ref = 'G26'; Value = '18' },
    @{ Ref = 'D27'; Value = '0.001219' },
    @{ Ref = 'E27'; Value = '0.16%' },
    @{ Ref = 'G27'; Value = '18' },
    @{ Ref = 'E28'; Value = '23.51%' },
    @{ Ref = 'G28'; Value = '18' },
    @{ Ref = 'G29'; Value = '18' },
    @{ Ref = 'G30'; Value = '18' },
    @{ Ref = 'G31'; Value = '18' },
    @{ Ref = 'G32'; Value = '18' },
    @{ Ref = 'G33'; Value = '18' },
    @{ Ref = 'G34'; Value = '18' },
    @{ Ref = 'G35'; Value = '18' },
    @{ Ref = 'G36'; Value = '18' },
    @{ Ref = 'G37'; Value = '18' },
    @{ Ref = 'G38'; Value = '18' },
    @{ Ref = 'G39'; Value = '18' },
    @{ Ref = 'D40'; Value = '0.03747' },
    @{ Ref = 'E40'; Value = '1.36%' },
    @{ Ref = 'G40'; Value = '18' },
    @{ Ref = 'D41'; Value = '0.005760' },
    @{ Ref = 'E41'; Value = '68.17%' },
    @{ Ref = 'G41'; Value = '18' },
    @{ Ref = 'E42'; Value = '1.43%' },
    @{ Ref = 'G42'; Value = '18' },
    @{ Ref = 'E43'; Value = '-18.01%' },
    @{ Ref = 'G43'; Value = '18' },
    @{ Ref = 'D44'; Value = '0.008277' },
    @{ Ref = 'E44'; Value = '-11.94%' },
    @{ Ref = 'G44'; Value = '18' },
    @{ Ref = 'E45'; Value = '0.20%' },
    @{ Ref = 'G45'; Value = '18' },
    @{ Ref = 'E46'; Value = '0.04%' },
    @{ Ref = 'G46'; Value = '18' },
    @{ Ref = 'D47'; Value = '0.05800' },
    @{ Ref = 'E47'; Value = '-51.65%' },
    @{ Ref = 'G47'; Value = '18' },
    @{ Ref = 'D48'; Value = '0.002573' },
    @{ Ref = 'G48'; Value = '18' },
    @{ Ref = 'E49'; Value = '0.04%' },
    @{ Ref = 'G49'; Value = '18' },
    @{ Ref = 'E50'; Value = '0.04%' },
    @{ Ref = 'G50'; Value = '18' },
    @{ Ref = 'G51'; Value = '18' }
)
foreach ($u in $textUpdates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
